$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2023-09-08 Friday" "2023-09-09 Saturday"

Replace-Text "11÷3=3, 2" "73÷7=10, 3"
Replace-Text "52÷6=8, 4" "23÷7=3, 2"
Replace-Text "44÷7=6, 2" "47÷4=11, 3"
Replace-Text "19÷2=9, 1" "13÷5=2, 3"
Replace-Text "89÷9=9, 8" "14÷5=2, 4"

Replace-Text "61÷4=15, 1" "82÷2=41, 0"
Replace-Text "66÷7=9, 3" "72÷3=24, 0"
Replace-Text "24÷6=4, 0" "38÷4=9, 2"
Replace-Text "63÷7=9, 0" "14÷5=2, 4"
Replace-Text "92÷3=30, 2" "14÷2=7, 0"

Replace-Text "20÷2=10, 0" "11÷6=1, 5"
Replace-Text "43÷4=10, 3" "54÷3=18, 0"
Replace-Text "83÷7=11, 6" "90÷9=10, 0"
Replace-Text "27÷9=3, 0" "94÷3=31, 1"
Replace-Text "67÷6=11, 1" "39÷4=9, 3"

Replace-Text "18÷7=2, 4" "32÷8=4, 0"
Replace-Text "29÷5=5, 4" "32÷7=4, 4"
Replace-Text "87÷8=10, 7" "11÷8=1, 3"
Replace-Text "70÷4=17, 2" "12÷7=1, 5"
Replace-Text "96÷3=32, 0" "15÷5=3, 0"

Replace-Text "78÷4=19, 2" "12÷6=2, 0"
Replace-Text "33÷3=11, 0" "10÷2=5, 0"
Replace-Text "81÷3=27, 0" "22÷6=3, 4"
Replace-Text "26÷3=8, 2" "96÷3=32, 0"
Replace-Text "53÷8=6, 5" "53÷4=13, 1"
